$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G ("Mã thuế" and everything after shifts right by one)
$ws.Columns("G:G").Insert()

# New header for the inserted column - use Formula so the inherited number format/style
# from the insert (matching column F) is preserved on the header cell
$ws.Range("G1").Formula = "Facebook"

# Data rows have nothing in the new Facebook column
$ws.Range("G2:G3").Clear()

# Match column F's width for the new column
$ws.Columns("G:G").ColumnWidth = $ws.Columns("F:F").ColumnWidth

# Update the active selection like in the authored workbook
$ws.Range("G2").Select()
